# Update cryptocurrency price/volume snapshot data (GitHub Actions refresh).
# Cell values that look like plain numbers (e.g. "570.33") are written with a
# leading apostrophe so Excel stores them as text, matching the source data's
# inlineStr / shared-string cells (e.g. "63.182.78" keeps its multi-dot format,
# "0.999" is not collapsed to 0.999 as a float/General-formatted number, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.182.78"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "2.560.16"
$ws.Range("E3").Value = "  +4.64%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'570.33"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "'148.09"
$ws.Range("E6").Value = "  +3.28%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -1.09%  "

$ws.Range("D9").Value = "2.561.29"
$ws.Range("E9").Value = "  +4.67%  "

$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").Value = "'5.60"
$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("D14").Value = "'27.34"
$ws.Range("E14").Value = "  +3.61%  "

$ws.Range("D15").Value = "3.019.49"
$ws.Range("E15").Value = "  +4.75%  "

$ws.Range("D16").Value = "63.140.51"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "'0.0000143"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "2.583.60"
$ws.Range("E18").Value = "  +5.86%  "

$ws.Range("E19").Value = "  +2.36%  "

$ws.Range("D20").Value = "'335.15"
$ws.Range("E20").Value = "  -2.03%  "

$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "'6.79"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "'64.93"
$ws.Range("E24").Value = "  -0.68%  "

$ws.Range("E25").Value = "  -2.79%  "

$ws.Range("E26").Value = "  +5.01%  "

$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").Value = "'1.51"
$ws.Range("E27").Value = "  +13.71%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.01"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Value = "'8.44"
$ws.Range("E29").Value = "  +2.68%  "

$ws.Range("D30").Value = "'7.28"
$ws.Range("E30").Value = "  +7.41%  "

$ws.Range("D31").Value = "0.0₃0821"
$ws.Range("E31").Value = "  +2.42%  "

$ws.Range("D32").Value = "'1.86"
$ws.Range("E32").Value = "  +1.44%  "

$ws.Range("D33").Value = "'177.60"
$ws.Range("E33").Value = "  +1.72%  "

$ws.Range("D34").Value = "'1.59"
$ws.Range("E34").Value = "  +6.60%  "

$ws.Range("D35").Value = "'416.03"
$ws.Range("E35").Value = "  +12.05%  "

$ws.Range("D36").Value = "'0.399"
$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").Value = "'4.40"
$ws.Range("E38").Value = "  -2.32%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  +2.35%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "'39.18"
$ws.Range("E42").Value = "  -2.37%  "

$ws.Range("D43").Value = "'152.39"
$ws.Range("E43").Value = "  +0.97%  "

$ws.Range("D44").Value = "'3.77"
$ws.Range("E44").Value = "  +1.16%  "

$ws.Range("D45").Value = "'20.85"
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("D46").Value = "'0.608"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("D47").Value = "'0.0963"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").Value = "'0.0522"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("D49").Value = "'0.0238"
$ws.Range("E49").Value = "  +5.64%  "

$ws.Range("D50").Value = "'18.48"
$ws.Range("E50").Value = "  +2.41%  "

$ws.Range("E51").Value = "  +3.14%  "
